$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "735. Asteroid Collision"
$ws.Range("B37").Value = "Medium"
$ws.Range("B37").Interior.Color = $ws.Range("B20").Interior.Color
$ws.Range("C37").Value = "Stacks"

$linkUrl = "https://leetcode.com/problems/asteroid-collision/solutions/3394436/c-java-python-javascript-stack-with-explaination/?envType=study-plan-v2&envId=leetcode-75 "
$ws.Hyperlinks.Add($ws.Range("E37"), $linkUrl) | Out-Null
$ws.Range("E37").Style = "Hyperlink"

$ws.Range("D37").Value = "The crux is to repeat the collision check at the top of the list. Just add to stack as you iterate the list, but consider the behavior for each case. Then at the end, populate the res list with the remaining stack elements."

$ws.Range("E41").Select() | Out-Null
